$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.555934307191982
$ws.Range("E2").Value = 2.196942688864458
$ws.Range("F2").Value = 2.833895656486887
$ws.Range("G2").Value = 3.429879838464347
$ws.Range("H2").Value = 3.965897297953021
$ws.Range("I2").Value = 4.434129702481918
$ws.Range("J2").Value = 4.833204475112849
$ws.Range("K2").Value = 5.164853783834539
$ws.Range("L2").Value = 5.431942562213596
$ws.Range("M2").Value = 5.629739127883019
$ws.Range("N2").Value = 5.762400891624407
$ws.Range("O2").Value = 5.83190203947099
$ws.Range("P2").Value = 5.837834837208962
$ws.Range("Q2").Value = 5.792535391542358
$ws.Range("R2").Value = 5.718709968800684
$ws.Range("S2").Value = 5.631227585998886
$ws.Range("T2").Value = 5.539622666095537
$ws.Range("U2").Value = 5.449811363270618
$ws.Range("V2").Value = 5.365271549340568
$ws.Range("W2").Value = 5.287855316535992
$ws.Range("X2").Value = 5.218348459567654
$ws.Range("Y2").Value = 5.156854923390787
$ws.Range("Z2").Value = 5.103059683262748
$ws.Range("AA2").Value = 5.056406965797862
$ws.Range("AB2").Value = 5.016219464211964
$ws.Range("AC2").Value = 4.981776487063499
$ws.Range("AD2").Value = 4.952363645454322
$ws.Range("AE2").Value = 4.927302961733269
$ws.Range("AF2").Value = 4.90969990275002
